# Remove the trailing duration annotation " (1 year 5 months)" that
# follows "May 2015 - Present" in the Globo.com experience entry.
#
# The paragraph currently contains three runs:
#   1) "May 2015 - Present"
#   2) " "                       (plain space)
#   3) "(1 year 5 months)"       (italic)
#
# We locate the text " (1 year 5 months)" (including its leading space)
# with Find and delete that exact Range. Deleting the found Range removes
# runs 2 and 3 in their entirety without touching/retyping run 1, so its
# original formatting/serialization (e.g. xml:space="preserve") is left
# completely untouched, matching the target diff.

$d = $word.ActiveDocument

$range = $d.Content
$found = $range.Find.Execute(
    " (1 year 5 months)",  # FindText
    $false,                 # MatchCase
    $false,                 # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    "",                     # ReplaceWith
    0                       # Replace (wdReplaceNone - search only)
)

if ($found) {
    $range.Delete()
}
